# Link the Keynote speaker "Julia Lane" entry to her new speaker page
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 is the "Keynote speaker" session; column F holds Speaker(s).
# Update the Julia Lane link to point at the new DIFA project speaker page.
$ws.Range("F4").Value = "[Julia Lane](https://dataifa.github.io/difa-project/julia_lane.html)"

# Update the active selection to match the author's final cursor position.
$ws.Range("F14").Select()
